$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Make room: insert 10 rows above the old "CCP (5 jours)" block
#    (old rows 19-23 -> become rows 29-33; Excel auto-shifts the
#    formulas that reference them).
# ------------------------------------------------------------------
$ws.Range("A19:A28").EntireRow.Insert()

# ------------------------------------------------------------------
# 2. The original bilateral block (rows 5-11) was simply "Bilat";
#    now that there are two liquidation periods it becomes
#    "Bilat (15 jours)".
# ------------------------------------------------------------------
$ws.Range("A5").Value = "Bilat (15 jours)"

# ------------------------------------------------------------------
# 3. Build the new "Bilat (5 jours)" block in rows 21-27, mirroring
#    the structure of the existing Bilat (15 jours) block (rows 5-11).
#    Copy the formatting first, then write labels/values/formulas.
# ------------------------------------------------------------------
$ws.Range("A5:L11").Copy()
$ws.Range("A21:L27").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A21").Value = "Bilat (5 jours)"

$ws.Range("B21").Value = "CVA"
$ws.Range("D21").Value = 0.034217
$ws.Range("E21").Value = 0.035531
$ws.Range("G21").Formula = "=AVERAGE(D21,E21)"
$ws.Range("I21").Formula = "=G21"

$ws.Range("B22").Value = "DVA"
$ws.Range("D22").Value = 0.032122
$ws.Range("E22").Value = 0.033329
$ws.Range("G22").Formula = "=AVERAGE(D22,E22)"
$ws.Range("I22").Formula = "=G22"

$ws.Range("B23").Value = "FVA"
$ws.Range("D23").Value = 0.000028
$ws.Range("E23").Value = 0.000405
$ws.Range("G23").Formula = "=AVERAGE(D23,E23)"
$ws.Range("I23").Formula = "=G23"

$ws.Range("B24").Value = "KVA_CCR_95"
$ws.Range("D24").Value = 0.385057
$ws.Range("E24").Value = 0.407967
$ws.Range("G24").Formula = "=AVERAGE(D24,E24)"
$ws.Range("I24").Formula = "=SUM(G24,G26,)"

$ws.Range("B25").Value = "KVA_CCR_99"
$ws.Range("D25").Value = 0.421126
$ws.Range("E25").Value = 0.444985
$ws.Range("G25").Formula = "=AVERAGE(D25,E25)"
$ws.Range("I25").Formula = "=SUM(G25,G27)"

$ws.Range("B26").Value = "KVA_CVA_95"
$ws.Range("D26").Value = 0.128822
$ws.Range("E26").Value = 0.136289
$ws.Range("G26").Formula = "=AVERAGE(D26,E26)"

$ws.Range("B27").Value = "KVA_CVA_99"
$ws.Range("D27").Value = 0.141843
$ws.Range("E27").Value = 0.149657
$ws.Range("G27").Formula = "=AVERAGE(D27,E27)"

# ------------------------------------------------------------------
# 4. The "Ratio bilat / CCP" column for the 5-day CCP block (now at
#    rows 29-33) must compare against the new Bilat (5 jours) block
#    (rows 21-25) instead of the 15-day one (rows 5-9).
# ------------------------------------------------------------------
$ws.Range("L29").Formula = "=I21/I29"
$ws.Range("L30").Formula = "=I22/I30"
$ws.Range("L31").Formula = "=I23/I31"
$ws.Range("L32").Formula = "=I24/I32"
$ws.Range("L33").Formula = "=I25/I33"

# ------------------------------------------------------------------
# 5. Cosmetic: restore the selection that was active when the author
#    saved the workbook.
# ------------------------------------------------------------------
$ws.Range("M20").Select()
